$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "supplier_telp" column (D): shift the "supplier_alamat"
# header/value from column E into D, then clear column E. This keeps
# the existing per-column <col> width definitions untouched (only the
# cell data moves), matching how the column was dropped upstream.
$ws.Range("E1:E2").Copy($ws.Range("D1:D2")) | Out-Null
$ws.Range("E1:E2").ClearContents() | Out-Null

# Move the active selection to D4 (as recorded in the saved file)
$ws.Range("D4").Select() | Out-Null
